$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# Target stored widths: col1=23, col2=25, col6=26.83203125 (engine snaps to 1/6 grid; closest achievable)
$ws.Columns.Item(1).ColumnWidth = 22.166666666666668
$ws.Columns.Item(2).ColumnWidth = 24.166666666666668
$ws.Columns.Item(6).ColumnWidth = 26.0

# --- Row 56: new section title ---
$ws.Range("A56").Value = "Updated directory and file names"
$ws.Range("A56").Font.Bold = $true

# --- Row 58: repeated header row ---
$ws.Range("A58").Value = "Project category directory"
$ws.Range("B58").Value = "Project directory"
$ws.Range("C58").Value = "OWL file"
$ws.Range("D58").Value = "Doc directory"
$ws.Range("E58").Value = "Conversion file"
$ws.Range("F58").Value = "Conversion Temporary file"
$ws.Range("A58:F58").Font.Bold = $true

# --- Rows 59-108: updated directory/file table ---
# row 59
$ws.Range("A59").Value = "Gates"
$ws.Range("B59").Value = "ELICIT"
$ws.Range("C59").Value = "gates_elicit.owl"
$ws.Range("D59").Value = "doc"
$ws.Range("E59").Value = "ELICIT_conversion.csv"
$ws.Range("F59").Value = "ELICIT_conversion_temp.csv"
$ws.Range("D59").Font.Color = 0
# row 60
$ws.Range("A60").Value = "Gates"
$ws.Range("B60").Value = "GAMIN"
$ws.Range("C60").Value = "gates_gamin.owl"
$ws.Range("D60").Value = "doc"
$ws.Range("E60").Value = "GAMIN_conversion.csv"
$ws.Range("F60").Value = "GAMIN_conversion_temp.csv"
$ws.Range("D60").Font.Color = 0
# row 61
$ws.Range("A61").Value = "Gates"
$ws.Range("B61").Value = "ganc"
$ws.Range("C61").Value = "gates_ganc.owl"
$ws.Range("D61").Value = "doc"
$ws.Range("E61").Value = "ganc_conversion.csv"
$ws.Range("F61").Value = "ganc_conversion_temp.csv"
$ws.Range("D61").Font.Color = 0
# row 62
$ws.Range("A62").Value = "Gates"
$ws.Range("B62").Value = "GEMS"
$ws.Range("C62").Value = "gates_gems.owl"
$ws.Range("D62").Value = "doc"
$ws.Range("E62").Value = "GEMS_conversion.csv"
$ws.Range("F62").Value = "GEMS_conversion_temp.csv"
$ws.Range("D62").Font.Color = 0
# row 63
$ws.Range("A63").Value = "Gates"
$ws.Range("B63").Value = "GEMS_HUAS"
$ws.Range("C63").Value = "gates_gems_huas.owl"
$ws.Range("D63").Value = "doc"
$ws.Range("E63").Value = "gems_huas_conversion.csv"
$ws.Range("F63").Value = "gems_huas_conversion_temp.csv"
$ws.Range("D63").Font.Color = 0
# row 64
$ws.Range("A64").Value = "Gates"
$ws.Range("B64").Value = "GEMS1A"
$ws.Range("C64").Value = "gates_gems1a.owl"
$ws.Range("D64").Value = "doc"
$ws.Range("E64").Value = "GEMS1a_conversion.csv"
$ws.Range("F64").Value = "GEMS1a_conversion_temp.csv"
$ws.Range("D64").Font.Color = 0
# row 65
$ws.Range("A65").Value = "Gates"
$ws.Range("B65").Value = "GEMS1A_HUAS"
$ws.Range("C65").Value = "gates_gems1a_huas.owl"
$ws.Range("D65").Value = "doc"
$ws.Range("E65").Value = "gems1a_huas_conversion.csv"
$ws.Range("F65").Value = "gems1a_huas_conversion_temp.csv"
$ws.Range("D65").Font.Color = 0
# row 66
$ws.Range("A66").Value = "Gates"
$ws.Range("B66").Value = "LLINEUP"
$ws.Range("C66").Value = "gates_llineup.owl"
$ws.Range("D66").Value = "doc"
$ws.Range("E66").Value = "LLINEUP_conversion.csv"
$ws.Range("F66").Value = "LLINEUP_conversion_temp.csv"
$ws.Range("D66").Font.Color = 0
# row 67
$ws.Range("A67").Value = "Gates"
$ws.Range("B67").Value = "MALED"
$ws.Range("C67").Value = "gates_maled.owl"
$ws.Range("D67").Value = "doc"
$ws.Range("E67").Value = "MALED_conversion.csv"
$ws.Range("D67").Font.Color = 0
# row 68
$ws.Range("A68").Value = "Gates"
$ws.Range("B68").Value = "PERCH"
$ws.Range("C68").Value = "gates_perch.owl"
$ws.Range("D68").Value = "doc"
$ws.Range("E68").Value = "PERCH_conversion.csv"
$ws.Range("F68").Value = "PERCH_conversion_temp.csv"
$ws.Range("D68").Font.Color = 0
# row 69
$ws.Range("A69").Value = "Gates"
$ws.Range("B69").Value = "PROVIDE"
$ws.Range("C69").Value = "gates_provide.owl"
$ws.Range("D69").Value = "doc"
$ws.Range("E69").Value = "PROVIDE_conversion.csv"
$ws.Range("F69").Value = "PROVIDE_conversion_temp.csv"
$ws.Range("D69").Font.Color = 0
# row 70
$ws.Range("A70").Value = "Gates"
$ws.Range("B70").Value = "SCORE_BURUNDI"
$ws.Range("C70").Value = "gates_score_burundi.owl"
$ws.Range("D70").Value = "doc"
$ws.Range("E70").Value = "SCORE_Burundi_conversion.csv"
$ws.Range("F70").Value = "SCORE_Burundi_conversion_temp.csv"
$ws.Range("D70").Font.Color = 0
# row 71
$ws.Range("A71").Value = "Gates"
$ws.Range("B71").Value = "SCORE_Five_Country"
$ws.Range("C71").Value = "gates_score_five_country.owl"
$ws.Range("D71").Value = "doc"
$ws.Range("E71").Value = "score_five_country_conversion.csv"
$ws.Range("F71").Value = "score_five_country_conversion_temp.csv"
$ws.Range("D71").Font.Color = 0
# row 72
$ws.Range("A72").Value = "Gates"
$ws.Range("B72").Value = "SCORE_MOZ"
$ws.Range("C72").Value = "gates_score_moz.owl"
$ws.Range("D72").Value = "doc"
$ws.Range("E72").Value = "score_moz_conversion.csv"
$ws.Range("F72").Value = "score_moz_conversion_temp.csv"
$ws.Range("D72").Font.Color = 0
# row 73
$ws.Range("A73").Value = "Gates"
$ws.Range("B73").Value = "SCORE_NIG"
$ws.Range("C73").Value = "gates_score_nig.owl"
$ws.Range("D73").Value = "doc"
$ws.Range("E73").Value = "score_nig_conversion.csv"
$ws.Range("F73").Value = "score_nig_conversion_temp.csv"
$ws.Range("D73").Font.Color = 0
# row 74
$ws.Range("A74").Value = "Gates"
$ws.Range("B74").Value = "SCORE_RWANDA"
$ws.Range("C74").Value = "gates_score_rwanda.owl"
$ws.Range("D74").Value = "doc"
$ws.Range("E74").Value = "SCORE_Rwanda_conversion.csv"
$ws.Range("F74").Value = "SCORE_Rwanda_conversion_temp.csv"
$ws.Range("D74").Font.Color = 0
# row 75
$ws.Range("A75").Value = "Gates"
$ws.Range("B75").Value = "SCORE_SEASONAL"
$ws.Range("C75").Value = "gates_score_seasonal.owl"
$ws.Range("D75").Value = "doc"
$ws.Range("E75").Value = "score_seasonal_conversion.csv"
$ws.Range("F75").Value = "score_seasonal_conversion_temp.csv"
$ws.Range("D75").Font.Color = 0
# row 76
$ws.Range("A76").Value = "Gates"
$ws.Range("B76").Value = "SCORE_sm_cohort"
$ws.Range("C76").Value = "gates_score_sm_cohort.owl"
$ws.Range("D76").Value = "doc"
$ws.Range("E76").Value = "SCORE_Sm_Cohort_conversion.csv"
$ws.Range("F76").Value = "SCORE_Sm_Cohort_conversion_temp.csv"
$ws.Range("D76").Font.Color = 0
# row 77
$ws.Range("A77").Value = "Gates"
$ws.Range("B77").Value = "score_sm_crt"
$ws.Range("C77").Value = "gates_score_sm_crt.owl"
$ws.Range("D77").Value = "doc"
$ws.Range("E77").Value = "score_sm_crt_conversion.csv"
$ws.Range("F77").Value = "score_sm_crt_conversion_temp.csv"
$ws.Range("D77").Font.Color = 0
# row 78
$ws.Range("A78").Value = "Gates"
$ws.Range("B78").Value = "Score_Zanzibar"
$ws.Range("C78").Value = "gates_score_zanzibar.owl"
$ws.Range("D78").Value = "doc"
$ws.Range("E78").Value = "score_zanzibar_conversion.csv"
$ws.Range("F78").Value = "score_zanzibar_conversion_temp.csv"
$ws.Range("D78").Font.Color = 0
# row 79
$ws.Range("A79").Value = "Gates"
$ws.Range("B79").Value = "SHINE"
$ws.Range("C79").Value = "gates_shine.owl"
$ws.Range("D79").Value = "doc"
$ws.Range("E79").Value = "SHINE_conversion.csv"
$ws.Range("F79").Value = "SHINE_conversion_temp.csv"
$ws.Range("D79").Font.Color = 0
# row 80
$ws.Range("A80").Value = "Gates"
$ws.Range("B80").Value = "SIP"
$ws.Range("C80").Value = "gates_sip.owl"
$ws.Range("D80").Value = "doc"
$ws.Range("E80").Value = "SIP_conversion.csv"
$ws.Range("F80").Value = "SIP_conversion_temp.csv"
$ws.Range("D80").Font.Color = 0
# row 81
$ws.Range("A81").Value = "Gates"
$ws.Range("B81").Value = "VIDA"
$ws.Range("C81").Value = "gates_vida.owl"
$ws.Range("D81").Value = "doc"
$ws.Range("E81").Value = "VIDA_conversion.csv"
$ws.Range("F81").Value = "VIDA_conversion_temp.csv"
$ws.Range("D81").Font.Color = 0
# row 82
$ws.Range("A82").Value = "Gates"
$ws.Range("B82").Value = "VIDA_HUCS_GAMBIA_MALI"
$ws.Range("C82").Value = "gates_vida_hucs_gambia_mali.owl"
$ws.Range("D82").Value = "doc"
$ws.Range("E82").Value = "VIDA_HUCS_GAMBIA_MALI_conversion.csv"
$ws.Range("F82").Value = "VIDA_HUCS_GAMBIA_MALI_conversion_temp.csv"
$ws.Range("D82").Font.Color = 0
# row 83
$ws.Range("A83").Value = "Gates"
$ws.Range("B83").Value = "VIDA_HUCS_KENYA"
$ws.Range("C83").Value = "gates_vida_hucs_kenya.owl"
$ws.Range("D83").Value = "doc"
$ws.Range("E83").Value = "VIDA_HUCS_KENYA_conversion.csv"
$ws.Range("F83").Value = "VIDA_HUCS_KENYA_conversion_temp.csv"
$ws.Range("D83").Font.Color = 0
# row 84
$ws.Range("A84").Value = "Gates"
$ws.Range("B84").Value = "washb_bangladesh"
$ws.Range("C84").Value = "gates_washb_bangladesh.owl"
$ws.Range("D84").Value = "doc"
$ws.Range("E84").Value = "WASHb_Bangladesh_conversion.csv"
$ws.Range("F84").Value = "WASHb_Bangladesh_conversion_temp.csv"
$ws.Range("D84").Font.Color = 0
# row 85
$ws.Range("A85").Value = "Gates"
$ws.Range("B85").Value = "washb_kenya"
$ws.Range("C85").Value = "gates_washb_kenya.owl"
$ws.Range("D85").Value = "doc"
$ws.Range("E85").Value = "WASHb_Kenya_conversion.csv"
$ws.Range("F85").Value = "WASHb_Kenya_conversion_temp.csv"
$ws.Range("D85").Font.Color = 0
# row 86
$ws.Range("A86").Value = "Gates"
$ws.Range("B86").Value = "WOMAN"
$ws.Range("C86").Value = "gates_woman.owl"
$ws.Range("D86").Value = "doc"
$ws.Range("E86").Value = "WOMAN_conversion.csv"
$ws.Range("F86").Value = "WOMAN_conversion_temp.csv"
$ws.Range("D86").Font.Color = 0
# row 87
$ws.Range("A87").Value = "General"
$ws.Range("B87").Value = "covid19_india"
$ws.Range("C87").Value = "general_covid19_india.owl"
$ws.Range("D87").Value = "doc"
$ws.Range("E87").Value = "covid19_india_conversion.csv"
$ws.Range("F87").Value = "covid19_india_conversion_temp.csv"
$ws.Range("D87").Font.Color = 0
# row 88
$ws.Range("A88").Value = "General"
$ws.Range("B88").Value = "hcq_non_randomized"
$ws.Range("C88").Value = "general_hcq_non_randomized.owl"
$ws.Range("D88").Value = "doc"
$ws.Range("E88").Value = "hcq_non_randomized_conversion.csv"
$ws.Range("F88").Value = "hcq_non_randomized_conversion_temp.csv"
$ws.Range("D88").Font.Color = 0
# row 89
$ws.Range("A89").Value = "General"
$ws.Range("B89").Value = "kalifabougou"
$ws.Range("C89").Value = "general_kalifabougou.owl"
$ws.Range("D89").Value = "doc"
$ws.Range("E89").Value = "kalifabougou_conversion.csv"
$ws.Range("F89").Value = "kalifabougou_conversion_temp.csv"
$ws.Range("D89").Font.Color = 0
# row 90
$ws.Range("A90").Value = "General"
$ws.Range("B90").Value = "nhs"
$ws.Range("C90").Value = "general_nhs.owl"
$ws.Range("D90").Value = "doc"
$ws.Range("E90").Value = "NHS_conversion.csv"
$ws.Range("F90").Value = "NHS_conversion_temp.csv"
$ws.Range("D90").Font.Color = 0
# row 91
$ws.Range("A91").Value = "General"
$ws.Range("B91").Value = "promote"
$ws.Range("C91").Value = "general_promote.owl"
$ws.Range("D91").Value = "doc"
$ws.Range("E91").Value = "promote_conversion.csv"
$ws.Range("F91").Value = "promote_conversion_temp.csv"
$ws.Range("D91").Font.Color = 0
# row 92
$ws.Range("A92").Value = "General"
$ws.Range("B92").Value = "umsp"
$ws.Range("C92").Value = "general_umsp.owl"
$ws.Range("D92").Value = "doc"
$ws.Range("E92").Value = "UMSP_conversion.csv"
$ws.Range("F92").Value = "UMSP_conversion_temp.csv"
$ws.Range("D92").Font.Color = 0
# row 93
$ws.Range("A93").Value = "ICEMR"
$ws.Range("B93").Value = "amazonia_Brazil"
$ws.Range("C93").Value = "icemr_amazonia_brazil.owl"
$ws.Range("D93").Value = "doc"
$ws.Range("E93").Value = "amazonia_brazil_conversion.csv"
$ws.Range("F93").Value = "amazonia_brazil_conversion_temp.csv"
$ws.Range("D93").Font.Color = 0
# row 94
$ws.Range("A94").Value = "ICEMR"
$ws.Range("B94").Value = "amazonia_Peru"
$ws.Range("C94").Value = "icemr_amazonia_peru.owl"
$ws.Range("D94").Value = "doc"
$ws.Range("E94").Value = "amazonia_peru_conversion.csv"
$ws.Range("F94").Value = "amazonia_peru_conversion_temp.csv"
$ws.Range("D94").Font.Color = 0
# row 95
$ws.Range("A95").Value = "ICEMR"
$ws.Range("B95").Value = "india_behavior"
$ws.Range("C95").Value = "icemr_india_behavior.owl"
$ws.Range("D95").Value = "doc"
$ws.Range("E95").Value = "India_behavior_conversion.csv"
$ws.Range("F95").Value = "India_behavior_conversion_temp.csv"
$ws.Range("D95").Font.Color = 0
# row 96
$ws.Range("A96").Value = "ICEMR"
$ws.Range("B96").Value = "india_cohort"
$ws.Range("C96").Value = "icemr_india_cohort.owl"
$ws.Range("D96").Value = "doc"
$ws.Range("E96").Value = "India_cohort_conversion.csv"
$ws.Range("F96").Value = "india_cohort_conversion_temp.csv"
$ws.Range("D96").Font.Color = 0
# row 97
$ws.Range("A97").Value = "ICEMR"
$ws.Range("B97").Value = "india_cx"
$ws.Range("C97").Value = "icemr_india_cx.owl"
$ws.Range("D97").Value = "doc"
$ws.Range("E97").Value = "India_cx_conversion.csv"
$ws.Range("F97").Value = "India_cx_conversion_temp.csv"
$ws.Range("D97").Font.Color = 0
# row 98
$ws.Range("A98").Value = "ICEMR"
$ws.Range("B98").Value = "india_fever_surv"
$ws.Range("C98").Value = "icemr_india_fever_surv.owl"
$ws.Range("D98").Value = "doc"
$ws.Range("E98").Value = "india_fever_surv_conversion.csv"
$ws.Range("F98").Value = "india_fever_surv_conversion_temp.csv"
$ws.Range("D98").Font.Color = 0
# row 99
$ws.Range("A99").Value = "ICEMR"
$ws.Range("B99").Value = "india_meghalaya"
$ws.Range("C99").Value = "icemr_india_meghalaya.owl"
$ws.Range("D99").Value = "doc"
$ws.Range("E99").Value = "india_meghalaya_conversion.csv"
$ws.Range("F99").Value = "india_meghalaya_conversion_temp.csv"
$ws.Range("D99").Font.Color = 0
# row 100
$ws.Range("A100").Value = "ICEMR"
$ws.Range("B100").Value = "india_severe_malaria"
$ws.Range("C100").Value = "icemr_india_severe_malaria.owl"
$ws.Range("D100").Value = "doc"
$ws.Range("E100").Value = "india_severe_malaria_conversion.csv"
$ws.Range("F100").Value = "india_severe_malaria_conversion_temp.csv"
$ws.Range("D100").Font.Color = 0
# row 101
$ws.Range("A101").Value = "ICEMR"
$ws.Range("B101").Value = "malawi"
$ws.Range("C101").Value = "icemr_malawi.owl"
$ws.Range("D101").Value = "doc"
$ws.Range("E101").Value = "malawi_conversion.csv"
$ws.Range("F101").Value = "malawi_conversion_temp.csv"
$ws.Range("D101").Font.Color = 0
# row 102
$ws.Range("A102").Value = "ICEMR"
$ws.Range("B102").Value = "prism"
$ws.Range("C102").Value = "icemr_prism.owl"
$ws.Range("D102").Value = "doc"
$ws.Range("E102").Value = "PRISM_conversion.csv"
$ws.Range("F102").Value = "PRISM_conversion_temp.csv"
$ws.Range("D102").Font.Color = 0
# row 103
$ws.Range("A103").Value = "ICEMR"
$ws.Range("B103").Value = "prism2"
$ws.Range("C103").Value = "icemr_prism2.owl"
$ws.Range("D103").Value = "doc"
$ws.Range("E103").Value = "PRISM2_conversion.csv"
$ws.Range("F103").Value = "PRISM2_conversion_temp.csv"
$ws.Range("D103").Font.Color = 0
# row 104
$ws.Range("A104").Value = "ICEMR"
$ws.Range("B104").Value = "prism2_border_cohort"
$ws.Range("C104").Value = "icemr_prism2_border_cohort.owl"
$ws.Range("D104").Value = "doc"
$ws.Range("E104").Value = "prism2_border_cohort_conversion.csv"
$ws.Range("F104").Value = "prism2_border_cohort_conversion_temp.csv"
$ws.Range("D104").Font.Color = 0
# row 105
$ws.Range("A105").Value = "ICEMR"
$ws.Range("B105").Value = "south_asia"
$ws.Range("C105").Value = "icemr_south_asia.owl"
$ws.Range("D105").Value = "doc"
$ws.Range("E105").Value = "south_asia_conversion.csv"
$ws.Range("F105").Value = "south_asia_conversion_temp.csv"
$ws.Range("D105").Font.Color = 0
# row 106
$ws.Range("A106").Value = "ICEMR"
$ws.Range("B106").Value = "southeast_asia"
$ws.Range("C106").Value = "icemr_southeast_asia.owl"
$ws.Range("D106").Value = "doc"
$ws.Range("E106").Value = "southeast_asia_conversion.csv"
$ws.Range("F106").Value = "southeast_asia_conversion_temp.csv"
$ws.Range("D106").Font.Color = 0
# row 107
$ws.Range("A107").Value = "ICEMR"
$ws.Range("B107").Value = "southern_africa"
$ws.Range("C107").Value = "icemr_southern_africa.owl"
$ws.Range("D107").Value = "doc"
$ws.Range("E107").Value = "southern_africa_conversion.csv"
$ws.Range("F107").Value = "southern_africa_conversion_temp.csv"
$ws.Range("D107").Font.Color = 0
# row 108
$ws.Range("A108").Value = "ICEMR"
$ws.Range("B108").Value = "west_africa"
$ws.Range("C108").Value = "icemr_west_africa.owl"
$ws.Range("D108").Value = "doc"
$ws.Range("E108").Value = "west_africa_conversion.csv"
$ws.Range("F108").Value = "west_africa_conversion_temp.csv"
$ws.Range("D108").Font.Color = 0

# --- Selection / view state (best effort) ---
$ws.Range("C55").Select()

# --- Leftover sortState artifact from prior sort operation on helper range ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("H59:H112"))
$sortObj.SortFields.Add($ws.Range("I59:I112"))
$sortObj.SetRange($ws.Range("H59:K112"))
$sortObj.Header = -4142
$sortObj.Apply()

